$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.104136666666667
$ws.Range("N2").Value = 3.31241
$ws.Range("O2").Value = 0.2772362398998524
$ws.Range("P2").Value = 0.2772362398998525
$ws.Range("Q2").Value = 29.85464717310777
$ws.Range("R2").Value = 268.69182455797
$ws.Range("S2").Value = 0.01967320539466541
$ws.Range("T2").Value = 0.01967320539466542
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.230986666666666
$ws.Range("N3").Value = 6.692959999999999
$ws.Range("O3").Value = 0.5601755411317187
$ws.Range("P3").Value = 0.5601755411317187
$ws.Range("Q3").Value = 60.3234380235911
$ws.Range("R3").Value = 542.9109422123199
$ws.Range("S3").Value = 0.03975111075569746
$ws.Range("T3").Value = 0.03975111075569747
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.647533
$ws.Range("N4").Value = 1.942599
$ws.Range("O4").Value = 0.1625882189684289
$ws.Range("P4").Value = 0.1625882189684289
$ws.Range("Q4").Value = 17.50858370305367
$ws.Range("R4").Value = 157.577253327483
$ws.Range("S4").Value = 0.01153756603997441
$ws.Range("T4").Value = 0.01153756603997441
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.104136666666667
$ws.Range("N5").Value = 3.31241
$ws.Range("O5").Value = 0.2772362398998524
$ws.Range("P5").Value = 0.2772362398998525
$ws.Range("Q5").Value = 381.5527306484633
$ws.Range("R5").Value = 3433.97457583617
$ws.Range("S5").Value = 0.2514303785075105
$ws.Range("T5").Value = 0.2514303785075106
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.230986666666666
$ws.Range("N6").Value = 6.692959999999999
$ws.Range("O6").Value = 0.5601755411317187
$ws.Range("P6").Value = 0.5601755411317187
$ws.Range("Q6").Value = 770.9544301946132
$ws.Range("R6").Value = 6938.589871751518
$ws.Range("S6").Value = 0.5080329627478567
$ws.Range("T6").Value = 0.5080329627478567
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.647533
$ws.Range("N7").Value = 1.942599
$ws.Range("O7").Value = 0.1625882189684289
$ws.Range("P7").Value = 0.1625882189684289
$ws.Range("Q7").Value = 223.765763599607
$ws.Range("R7").Value = 2013.891872396463
$ws.Range("S7").Value = 0.1474540898796681
$ws.Range("T7").Value = 0.1474540898796682
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.104136666666667
$ws.Range("N8").Value = 3.31241
$ws.Range("O8").Value = 0.2772362398998524
$ws.Range("P8").Value = 0.2772362398998525
$ws.Range("Q8").Value = 9.306479415617778
$ws.Range("R8").Value = 83.75831474056001
$ws.Range("S8").Value = 0.00613265599767647
$ws.Range("T8").Value = 0.006132655997676473
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.230986666666666
$ws.Range("N9").Value = 6.692959999999999
$ws.Range("O9").Value = 0.5601755411317187
$ws.Range("P9").Value = 0.5601755411317187
$ws.Range("Q9").Value = 18.80440358215111
$ws.Range("R9").Value = 169.23963223936
$ws.Range("S9").Value = 0.0123914676281646
$ws.Range("T9").Value = 0.01239146762816461
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.647533
$ws.Range("N10").Value = 1.942599
$ws.Range("O10").Value = 0.1625882189684289
$ws.Range("P10").Value = 0.1625882189684289
$ws.Range("Q10").Value = 5.457886435042668
$ws.Range("R10").Value = 49.120977915384
$ws.Range("S10").Value = 0.003596563048786326
$ws.Range("T10").Value = 0.003596563048786328
